# Updates the weekly price data for Hortaliza / Espárragos (Vega Monumental
# Concepción). The underlying edit re-shuffles the per-row data (date,
# volume, min/max/weighted prices, variety, quality and origin) across
# rows 2-19 while leaving row 8 untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new value, matching the target OOXML.
$changes = @{
    2  = @{ D = 44477; J = 500;  K = 1400; L = 1500; M = 1460;                          P = 1460 }
    3  = @{ D = 44860; J = 1100; K = 1500; L = 1700; M = 1609;                          P = 1609 }
    4  = @{ D = 44875; J = 300;  K = 1500; L = 1600; M = 1550;                          P = 1550 }
    5  = @{ D = 44519; J = 250;                      M = 1240; O = 'Provincia de Linares'; P = 1240 }
    6  = @{ D = 44511; I = 'Primera'; J = 600; K = 1300; L = 1400; M = 1350; O = 'Provincia de Linares'; P = 1350 }
    7  = @{ D = 44876; H = 'Sin especificar'; J = 350; K = 1500; L = 1600; M = 1557;    P = 1557 }
    9  = @{ D = 44510; J = 600;  K = 1300; L = 1400; M = 1350;                          P = 1350 }
    10 = @{ D = 44489; J = 600;  K = 1400; L = 1500; M = 1450;                          P = 1450 }
    11 = @{ D = 44839; J = 500;  K = 1700; L = 1800; M = 1760;                          P = 1760 }
    12 = @{ D = 44526; J = 100;            L = 1600; M = 1550;                          P = 1550 }
    13 = @{ D = 44468; H = 'Verde'; J = 500; K = 1800; L = 2000; M = 1920;              P = 1920 }
    14 = @{ D = 44524; J = 200;  K = 1500; L = 1600; M = 1550; O = 'Provincia de Talca'; P = 1550 }
    15 = @{ D = 44868; J = 1000; K = 1200; L = 1300; M = 1250; O = 'Región del Maule';  P = 1250 }
    16 = @{ D = 44868; I = 'Segunda';      K = 1000; L = 1000; M = 1000; O = 'Región del Maule'; P = 1000 }
    17 = @{ D = 44545; I = 'Primera'; J = 550; K = 1700; L = 1800; M = 1755;            P = 1755 }
    18 = @{ D = 44881; J = 200;  K = 2600; L = 2700; M = 2650;                          P = 2650 }
    19 = @{ D = 44881; I = 'Segunda'; J = 100; K = 2400; L = 2400; M = 2400;            P = 2400 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
